# "working sript for multi substrat mm"
#
# 1) Rename the first two data sheets to the underscore/ASCII-only names.
# 2) Add a blank (empty-string) label cell at A2 on the "PD (5 mM NAD)"
#    sheet so its second header row lines up with the other sheets.
# 3) Clean up AZ13 on that same sheet: it held "0.378" with a trailing
#    stray carriage return (imported as literal "_x000d_") - strip it so
#    the cell just reads 0.378.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("NAD (500 mM PD)")
$ws1.Name = "NAD_PD_500nM"

$ws2 = $wb.Worksheets.Item("PD (5 mM NAD)")
$ws2.Name = "PD_NAD_5nM"

# Write A2 as text (a leading apostrophe forces text interpretation even
# though the content is empty), then drop the formatting it picked up so
# the cell keeps the workbook's default (unstyled) look.
$a2 = $ws2.Range("A2")
$a2.Value = "'"
$a2.ClearFormats()

# Same trick for AZ13: force text so "0.378" isn't reinterpreted as a
# number, then strip the formatting the coercion added.
$az13 = $ws2.Range("AZ13")
$az13.NumberFormat = "@"
$az13.Value = "0.378"
$az13.ClearFormats()
